$d = $word.ActiveDocument

# The document currently ends with a run of identical paragraphs:
#   "Documento de Evidencias - DemoBlaze" (sz=36)
# Append three more paragraphs with the same text/formatting at the end
# of the body (contemplating invalid credentials scenarios).

for ($i = 0; $i -lt 3; $i++) {
    $lastPara = $d.Paragraphs.Last
    $r = $lastPara.Range
    $r.InsertParagraphAfter()
    $d.Paragraphs.Last.Range.Text = "Documento de Evidencias - DemoBlaze"
}

Write-Host "Paragraphs count: " $d.Paragraphs.Count
